$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - copy style from existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-46
$values = @(
    @(5, 6),
    @(7, 8),
    @(5, 5),
    @(8, 9),
    @(10, 10),
    @(4, 5),
    @(5, 6),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(4, 5),
    @(5, 5),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(5, 5),
    @(5, 6),
    @(6, 6),
    @(7, 8),
    @(11, 11),
    @(2, 3),
    @(8, 8),
    @(8, 9),
    @(4, 4),
    @(8, 8),
    @(5, 5),
    @(8, 9),
    @(6, 6),
    @(6, 7),
    @(9, 9),
    @(7, 8),
    @(7, 7),
    @(7, 8),
    @(6, 8),
    @(7, 8),
    @(5, 7),
    @(7, 7),
    @(5, 6),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
